$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows (Uzbek, Vietnamese) which drop out of the data set.
$ws.Rows("22:23").Delete()

# Sorted (descending by 1994 value) language list that replaces rows 2-21.
$data = @(
    @("English", 26.91284187180431),
    @("Spanish", 8.232501437222416),
    @("Japanese", 7.939359535690782),
    @("Chinese", 6.608605746688448),
    @("German", 6.473975829475416),
    @("Arabic", 4.795307166969188),
    @("Portuguese", 3.889094560556003),
    @("Russian", 3.7749682316215),
    @("French", 3.721618929819467),
    @("Italian", 3.609502166086786),
    @("Malay-Indonesian", 2.71116628867311),
    @("Dutch", 1.663936414071622),
    @("Persian", 1.44589585184675),
    @("Korean", 1.345097265967378),
    @("Turkish", 1.314337776808709),
    @("Thai", 1.011391715893257),
    @("Polish", 0.801400365784672),
    @("Urdu", 0.7855456081773862),
    @("Swedish", 0.5198668517948949),
    @("Bengali", 0.4031970813447354)
)

$row = 2
foreach ($pair in $data) {
    $ws.Cells.Item($row, 1).Value = $pair[0]
    $ws.Cells.Item($row, 2).Value = $pair[1]
    $row = $row + 1
}
